# Gantt Tasks.xlsx edit script
# Reworks the task table: new task breakdown (UCD, testing tasks, documentation, ...),
# updated durations/dependencies, row 22-27 removed, some column widths changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152

# --- Row 1 (header) ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Task"
$ws.Range("C1").Value = "Duration (In Days)"
$ws.Range("D1").Value = "Dependencies"

# --- Row 2 ---
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "Analyse UCD"
$ws.Range("C2").Value = 1

# --- Row 3 ---
$ws.Range("A3").Value = "B"
$ws.Range("B3").Value = "Analyse Domain Model"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "A"
$ws.Range("D3").HorizontalAlignment = $xlRight

# --- Row 4 ---
$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = "Analyse ERD"
$ws.Range("C4").Value = 1
$ws.Range("D4").Clear()
$ws.Rows.Item(4).RowHeight = 15.75

# --- Row 5 ---
$ws.Range("A5").Value = "D"
$ws.Range("B5").Value = "Webscraper Schrijven"
$ws.Range("C5").Value = 2

# --- Row 6 ---
$ws.Range("A6").Value = "E"
$ws.Range("B6").Value = "Webscraper Testen "
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "D"
$ws.Range("D6").HorizontalAlignment = $xlRight

# --- Row 7 ---
$ws.Range("A7").Value = "F"
$ws.Range("B7").Value = "NBB Scrapen Via API"
$ws.Range("C7").Value = 2

# --- Row 8 ---
$ws.Range("A8").Value = "G"
$ws.Range("B8").Value = "NBB Scrapen Testen"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "F"

# --- Row 9 ---
$ws.Range("A9").Value = "H"
$ws.Range("B9").Value = "Data Ophalen uit kmos.csv"
$ws.Range("C9").Value = 1
$ws.Range("D9").Clear()

# --- Row 10 ---
$ws.Range("A10").Value = "I"
$ws.Range("B10").Value = "Data Ophalen testen"
$ws.Range("C10").Value = 1
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "H"

# --- Row 11 ---
$ws.Range("A11").Value = "J"
$ws.Range("B11").Value = "Webscraper, NBB en csv combineren in 1 script"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "E,G,I"

# --- Row 12 ---
$ws.Range("A12").Value = "K"
$ws.Range("B12").Value = "Opzetten server met databank"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "C"

# --- Row 13 ---
$ws.Range("A13").Value = "L"
$ws.Range("B13").Value = "Script koppelen met databank (opslaan)"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "J,K"

# --- Row 14 ---
$ws.Range("A14").Value = "M"
$ws.Range("B14").Value = "Beginnen ophalen van data"
$ws.Range("C14").Value = "x"
$ws.Range("C14").HorizontalAlignment = $xlRight
$ws.Range("D14").Value = "L"

# --- Row 15 ---
$ws.Range("A15").Value = "N "
$ws.Range("B15").Value = "Full-Text-Search gebruiken voor de score"
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "M"

# --- Row 16 ---
$ws.Range("A16").Value = "O"
$ws.Range("B16").Value = "Machine Learning Model schrijven"
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = "N"
$ws.Range("D16").HorizontalAlignment = $xlRight

# --- Row 17 ---
$ws.Range("A17").Value = "P"
$ws.Range("B17").Value = "Backend vertalen naar python"
$ws.Range("C17").Value = 20

# --- Row 18 ---
$ws.Range("A18").Value = "Q"
$ws.Range("B18").Value = "Frontend-koppelen met nieuwe API"
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = "P"
$ws.Range("D18").HorizontalAlignment = $xlRight

# --- Row 19 ---
$ws.Range("A19").Value = "R"
$ws.Range("B19").Value = "Frontend Optimaliseren"
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = "Q"
$ws.Range("D19").HorizontalAlignment = $xlRight

# --- Row 20 ---
$ws.Range("A20").Value = "S"
$ws.Range("B20").Value = "Insights Data"
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = "M,P"
$ws.Range("D20").HorizontalAlignment = $xlRight

# --- Row 21 ---
$ws.Range("A21").Value = "T"
$ws.Range("B21").Value = "Documentatie"
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = "A, B, C, D, E, F, G, H, I, J, K, L, M, N, O, P ,Q, R,S"
$ws.Range("D21").HorizontalAlignment = $xlRight

# --- Remove old leftover rows 22-27 (task letters U-Z, now unused) ---
$ws.Range("A22:A27").EntireRow.Delete()

# --- Column width tweaks ---
$ws.Columns.Item(4).ColumnWidth = 40.584
$ws.Columns.Item(8).ColumnWidth = 7.1666666666666667
$ws.Columns.Item(9).ColumnWidth = 42.417
$ws.Columns.Item(10).ColumnWidth = 16.1666666666666667
$ws.Columns.Item(11).ColumnWidth = 12.917

# --- Selection state as left by the author ---
$ws.Range("D26").Select() | Out-Null
